$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the bold/centered/bordered header style from the quarter-label
# cells A2:A66, leaving the header row (A1:B1) formatting untouched.
$ws.Range("A2:A66").Style = "Normal"
